{"js": "// Word Use Case Specification edit:\n//  1. Remove the stray \"_GoBack\" bookmark that originally sat right after \"UC009\".\n//  2. Tighten the table cross-reference from \" (table 3.1)\" to \" (table 3)\" and\n//     re-plant the \"_GoBack\" bookmark between the \"3\" and the closing \")\".\n//  3. Tighten the caption text from \"3.1-Output\" to \"3-Output\", keeping it split\n//     across two runs (\"3\" and \"-Output\") the way the original author's edit left it.\n\n// --- 1. Drop the old \"_GoBack\" bookmark (previously right after \"UC009\"). ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. \" (table 3.1)\" -> \" (table 3)\" with \"_GoBack\" re-inserted before the \")\". ---\nconst tableRefResults = context.document.body.search(\" (table 3.1)\", { matchCase: true });\nawait context.sync();\n\nif (tableRefResults.items.length > 0) {\n  const tableRefRange = tableRefResults.items[0];\n  // Split so the boundary lands right after the \"3\" (keeps the leading \" (table 3\").\n  const tableRefParts = tableRefRange.split([\"3\"], false);\n  tableRefParts.load(\"text\");\n  await context.sync();\n\n  const beforeParen = tableRefParts.items[0]; // \" (table 3\"\n  const afterDigits = tableRefParts.items[1]; // \".1)\"\n\n  // Plant the bookmark right at the \"3\" | \".1)\" boundary.\n  beforeParen.getRange(\"After\").insertBookmark(\"_GoBack\");\n  await context.sync();\n\n  // Remove the \".1\" left-over so only \")\" remains, in its own run.\n  const leftoverResults = afterDigits.search(\".1\", { matchCase: true });\n  await context.sync();\n  if (leftoverResults.items.length > 0) {\n    leftoverResults.items[0].delete();\n    await context.sync();\n  }\n}\n\n// --- 3. \"3.1-Output\" -> \"3-Output\", preserved as two sibling runs (\"3\" + \"-Output\"). ---\nconst captionResults = context.document.body.search(\"3.1-Output\", { matchCase: true });\nawait context.sync();\n\nif (captionResults.items.length > 0) {\n  const captionRange = captionResults.items[0];\n  const captionParts = captionRange.split([\"3\"], false);\n  captionParts.load(\"text\");\n  await context.sync();\n\n  const restPart = captionParts.items[1]; // \".1-Output\"\n\n  // Drop the \".1\" so the run reads \"-Output\".\n  const dotResults = restPart.search(\".1\", { matchCase: true });\n  await context.sync();\n  if (dotResults.items.length > 0) {\n    dotResults.items[0].delete();\n    await context.sync();\n  }\n\n  // The text edit above lets the engine re-merge \"3\" and \"-Output\" into one run\n  // (identical run formatting). Force the boundary to persist as two separate\n  // runs by planting a temporary bookmark right at the seam, then removing it\n  // again - this mirrors how the original edit ended up split across two runs.\n  const rejoined = context.document.body.search(\"3-Output\", { matchCase: true });\n  await context.sync();\n  if (rejoined.items.length > 0) {\n    const rejoinedParts = rejoined.items[0].split([\"3\"], false);\n    rejoinedParts.load(\"text\");\n    await context.sync();\n    rejoinedParts.items[0].getRange(\"After\").insertBookmark(\"__tmp_split_seam\");\n    await context.sync();\n    context.document.deleteBookmark(\"__tmp_split_seam\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word Use Case Specification edit:\n#  1. Remove the stray \"_GoBack\" bookmark that originally sat right after \"UC009\".\n#  2. Tighten the table cross-reference from \" (table 3.1)\" to \" (table 3)\" and\n#     re-plant the \"_GoBack\" bookmark between the \"3\" and the closing \")\".\n#  3. Tighten the caption text from \"3.1-Output\" to \"3-Output\", keeping it split\n#     across two runs (\"3\" and \"-Output\") the way the original author's edit left it.\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the old \"_GoBack\" bookmark (previously right after \"UC009\"). ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2. \" (table 3.1)\" -> \" (table 3)\" with \"_GoBack\" re-inserted before the \")\". ---\n$tableRefRange = $d.Content\n$tableRefFound = $tableRefRange.Find.Execute(\" (table 3.1)\")\nif ($tableRefFound) {\n    $matchStart = $tableRefRange.Start\n    $matchEnd = $tableRefRange.End\n\n    # \" (table 3\" is the first 9 characters of the match; the boundary sits\n    # right after the \"3\" and before \".1)\".\n    $boundary = $matchStart + 9\n\n    # Plant the bookmark at that boundary - this also forces the run to split\n    # into \" (table 3\" and \")\" instead of staying a single run.\n    $boundaryRange = $d.Range($boundary, $boundary)\n    $d.Bookmarks.Add(\"_GoBack\", $boundaryRange)\n\n    # Remove the leftover \".1\" so only \")\" remains after the bookmark.\n    $leftoverRange = $d.Range($boundary, $boundary + 2)\n    $leftoverRange.Delete()\n}\n\n# --- 3. \"3.1-Output\" -> \"3-Output\", preserved as two sibling runs (\"3\" + \"-Output\"). ---\n$captionRange = $d.Content\n$captionFound = $captionRange.Find.Execute(\"3.1-Output\")\nif ($captionFound) {\n    $capStart = $captionRange.Start\n\n    # Drop the \".1\" right after the \"3\".\n    $dotRange = $d.Range($capStart + 1, $capStart + 3)\n    $dotRange.Delete()\n\n    # Deleting the text lets Word's writer re-merge \"3\" and \"-Output\" into a\n    # single run (identical run formatting). Force the boundary to persist as\n    # two separate runs by planting a temporary bookmark right at the seam,\n    # then removing it again - mirrors how the original edit left two runs.\n    $seamScan = $d.Content\n    $seamFound = $seamScan.Find.Execute(\"3-Output\")\n    if ($seamFound) {\n        $seam = $seamScan.Start + 1\n        $seamRange = $d.Range($seam, $seam)\n        $d.Bookmarks.Add(\"__tmp_split_seam\", $seamRange)\n        $d.Bookmarks(\"__tmp_split_seam\").Delete()\n    }\n}\n"}
